$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.314.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.114.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.48%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.04%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.107.29"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.51%  "
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +14.09%  "
$ws.Range("E11").Value = "  +7.68%  "
$ws.Range("E12").Value = "  +4.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.37%  "
$ws.Range("E14").Value = "  +6.06%  "
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.632.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.224.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.114.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.16%  "
$ws.Range("E21").Value = "  +4.42%  "
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("E23").Value = "  +7.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.78%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("E29").Value = "  +5.23%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.04%  "
$ws.Range("E33").Value = "  +4.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0868"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +16.77%  "
$ws.Range("E36").Value = "  +7.12%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.35%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +20.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "440.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.924.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.02%  "
$ws.Range("E43").Value = "  +5.43%  "
$ws.Range("E44").Value = "  +12.03%  "
$ws.Range("E45").Value = "  +5.21%  "
$ws.Range("E46").Value = "  +8.81%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.93%  "
